$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: clear the "$policy: Policy" text from B7 (keep its existing style)
$ws.Range("B7").Value = ""

# Row 8: B8 becomes a two-line Java snippet that now declares + assigns the
# policy variable; C8/D8 drop the leading "$policy." -> "policy."
$ws.Range("B8").Value = "Policy policy = new Policy();`npolicy.setId(`"`$param`");"
$ws.Range("B8").WrapText = $true
$ws.Rows("8").RowHeight = 48

$ws.Range("C8").Value = "policy.setName(`"`$param`");"
$ws.Range("D8").Value = "policy.setAmount(`$param);"

# Move the active selection from C12 to D9
$ws.Range("D9").Select()
